$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Student"
$ws.Range("A2").Value = "b"
$ws.Range("A3").Value = "math"
$ws.Range("B3").Value = 90
$ws.Range("A5").Value = "Student Average:"
$ws.Range("B5").Value = 90
$ws.Range("A6").Value = "Average:"
$ws.Range("B6").Value = 90
